{"js": "// Update the signatory's name for the first outorgante (\"Manuel Dias\" ->\n// \"Manuel In\u00e1cio Veladas Dias\") in the protocol document body.\nconst body = context.document.body;\nconst results = body.search(\"Manuel Dias\", { matchCase: true, matchWholeWord: false });\nresults.load(\"text\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  // Replace text in place so the existing run formatting (Calibri rPr) is kept.\n  results.items[0].insertText(\"Manuel In\u00e1cio Veladas Dias\", Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Update the signatory's name for the first outorgante (\"Manuel Dias\" ->\n# \"Manuel In\u00e1cio Veladas Dias\") in the protocol document body, while\n# leaving every other run in the paragraph untouched.\n\n$d = $word.ActiveDocument\n\n$searchText  = \"Manuel Dias\"\n$replaceText = \"Manuel In\u00e1cio Veladas Dias\"\n\n# Locate the text to replace.\n$target = $d.Content\n$found = $target.Find.Execute($searchText)\nif (-not $found) {\n    throw \"Text '$searchText' not found\"\n}\n\n$targetStart = $target.Start\n$targetEnd   = $target.End\n\n# Pin the run that immediately follows the match with a throwaway bookmark\n# so that writing the new text into the matched run doesn't fold the\n# following run(s) into it.\n$docEnd = $d.Content.End\n$pinEnd = [Math]::Min($targetEnd + 2, $docEnd)\n$pinBookmarkAdded = $false\nif ($pinEnd -gt $targetEnd) {\n    $d.Bookmarks.Add(\"__pin_next\", $d.Range($targetEnd, $pinEnd)) | Out-Null\n    $pinBookmarkAdded = $true\n}\n\n# Scope the edit itself to a throwaway bookmark around just the matched\n# run so the replacement doesn't disturb the runs on either side.\n$d.Bookmarks.Add(\"__tmp_edit\", $d.Range($targetStart, $targetEnd)) | Out-Null\n$editRange = $d.Bookmarks.Item(\"__tmp_edit\").Range\n$editRange.Text = $replaceText\n$d.Bookmarks.Item(\"__tmp_edit\").Delete()\n\nif ($pinBookmarkAdded) {\n    $d.Bookmarks.Item(\"__pin_next\").Delete()\n}\n"}
